$d = $word.ActiveDocument
$newText = "Dates à utiliser pour la Campagne Leo: 14-23 avril, 14-23 mai"
$marker = "Dates à utiliser pour la Campagne"

$count = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*$marker*") {
        $count = $count + 1

        # Remove the paragraph's existing text (but keep the paragraph mark),
        # then insert a brand-new, unformatted run with the replacement text.
        $rng = $p.Range
        $rng.End = $rng.End - 1
        $rng.Delete()

        $rng2 = $p.Range
        $rng2.End = $rng2.End - 1
        $rng2.InsertAfter($newText)
    }
}

Write-Host ("Replaced " + $count + " paragraph(s)")
